$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 29. This shifts existing rows 29:118 down to 30:119,
# carrying their formatting (including the date style on column D) along with them.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record's data.
# (Same category/product metadata as the rest of the sheet, with the specific
# changed fields from the diff: D, N, O, P, Q, S, T.)
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(29, 3).Value = "La Araucanía"
$ws.Cells.Item(29, 4).Value = 45148
$ws.Cells.Item(29, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100108
$ws.Cells.Item(29, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(29, 9).Value = 100108004
$ws.Cells.Item(29, 10).Value = "Papaya"
$ws.Cells.Item(29, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 80
$ws.Cells.Item(29, 14).Value = 25000
$ws.Cells.Item(29, 15).Value = 25000
$ws.Cells.Item(29, 16).Value = 25000
$ws.Cells.Item(29, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(29, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(29, 19).Value = 2500
$ws.Cells.Item(29, 20).Value = 10
